$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.037010908126831
$ws.Range("B1").Value = 3.403537511825562
$ws.Range("C1").Value = 2.872161626815796
$ws.Range("D1").Value = 2.122005462646484
$ws.Range("E1").Value = 1.233277559280396
